$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 06 (2022-12-06, row 76): end work - fill in HORA F (C) and PAUSAS (E)
$ws.Range("C76").Value = 0.5416666666666666
$ws.Range("E76").Value = 0.08263888888888889

# Day 07 (2022-12-07, row 77): full day worked - HORA I, HORA F, PAUSAS, ASSUNTO, PRODUÇÃO
$ws.Range("B77").Value = 0.35625
$ws.Range("C77").Value = 0.7291666666666666
$ws.Range("E77").Value = 0.09861111111111111
$ws.Range("G77").Value = "Estágio + HARD"
$ws.Range("H77").Value = "Estágio + Hard + Atividade voluntária no Alpha EdTech "

# Day 08 (2022-12-08, row 78): init work - HORA I, ASSUNTO, PRODUÇÃO
$ws.Range("B78").Value = 0.3541666666666667
$ws.Range("G78").Value = "Estágio"
$ws.Range("H78").Value = "Estágio"

$ws.Range("H78").Select()
